$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (row 9), previously empty
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail"; change to "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row; delete it entirely,
# shifting all subsequent rows up by one.
$ws.Rows.Item(11).Delete()

# After the deletion, the row that used to be "Case Sensitive" (old row 15) is now row 14.
# Set its previously-empty Value cell to the text "true" (not the Boolean TRUE).
# A leading apostrophe forces text entry; then re-apply the original cell
# formatting (copied from a neighboring data cell) so the style index is
# unaffected by the quote-prefix that Value-assignment would otherwise add.
$ws.Range("B14").Value = "'true"
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
